# Applies the "Converted to S2O" change-log update + journal publishing
# model correction described in the commit.
#
# Summary of the edit:
#  - "Journal List" sheet: 5 journals (rows 2-6) that were "Subscribe to
#    Open"/"Sponsored S2O" are changed to "Hybrid" in the Publishing Model
#    column (F), reflecting that they converted away from S2O.
#  - "Journal Changes Notes" sheet: the previously-blank log rows 6-10 are
#    filled in with a copy of the (pre-change) journal details plus the
#    change type "Publishing Model" and change note "Converted to S2O".
#  - Various view-only changes (active sheet/cell selection, column widths).

$wb = $excel.ActiveWorkbook

$wsList   = $wb.Worksheets.Item("Journal List")
$wsNotes  = $wb.Worksheets.Item("Journal Changes Notes")
$wsTerms  = $wb.Worksheets.Item("Journal Changes Terms")

# ---------------------------------------------------------------------
# 1. Snapshot the current ("before") values of the 5 affected journal
#    rows on "Journal List" (columns A-H) before we overwrite column F.
# ---------------------------------------------------------------------
$journalRows = 2..6
$snapshot = @{}
foreach ($r in $journalRows) {
    $rowVals = @()
    for ($c = 1; $c -le 8; $c++) {
        $rowVals += $wsList.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# ---------------------------------------------------------------------
# 2. Update "Journal List" Publishing Model column (F) to "Hybrid" for
#    the 5 journals that converted away from Subscribe to Open.
# ---------------------------------------------------------------------
foreach ($r in $journalRows) {
    $wsList.Range("F$r").Value = "Hybrid"
}

# ---------------------------------------------------------------------
# 3. Log the change on "Journal Changes Notes" (rows 6-10), copying the
#    pre-change journal details and recording the change type/note.
# ---------------------------------------------------------------------
$logRow = 6
foreach ($r in $journalRows) {
    $vals = $snapshot[$r]
    for ($c = 1; $c -le 8; $c++) {
        $wsNotes.Cells.Item($logRow, $c).Value = $vals[$c - 1]
    }
    $wsNotes.Cells.Item($logRow, 9).Value  = "Publishing Model"
    $wsNotes.Cells.Item($logRow, 10).Value = "Converted to S2O"
    $logRow++
}

# ---------------------------------------------------------------------
# 4. Cosmetic view changes: column widths + active sheet/cell selection.
# ---------------------------------------------------------------------
$wsList.Columns.Item(1).ColumnWidth = 40.88671875
$wsNotes.Columns.Item(5).ColumnWidth = 31.21875
$wsNotes.Columns.Item(6).ColumnWidth = 23.33203125

$wsList.Range("E9").Select()
$wsNotes.Range("F16").Select()
$wsTerms.Range("A10").Select()

$wsNotes.Activate()
